$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; temporarily unprotect so the refreshed figures
# can be written, then restore protection once the edits are in place.
$ws.Unprotect()

# Update the disclosure text: model holdings "as of" date moves from 2021-03-29 to 2021-03-30
$ws.Range("A37").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-30 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for each holding row
$ws.Range("D2").Value = 0.03824687301751194
$ws.Range("E2").Value = 0.005080109417741285
$ws.Range("D3").Value = 0.02158618783985694
$ws.Range("E3").Value = 0.002356637863314859
$ws.Range("D4").Value = 0.01988853701663053
$ws.Range("E4").Value = 0.004984423676012373
$ws.Range("D5").Value = 0.04047056950337148
$ws.Range("E5").Value = -0.006634078212290562
$ws.Range("D6").Value = 0.03725157657988404
$ws.Range("E6").Value = 0.003532182103610726
$ws.Range("D7").Value = 0.02104128301357739
$ws.Range("E7").Value = -0.002130956993413369
$ws.Range("D8").Value = 0.03787485291618704
$ws.Range("E8").Value = -0.005560498220640531
$ws.Range("D9").Value = 0.02126596250209122
$ws.Range("E9").Value = 0.0004600239212437796
$ws.Range("D10").Value = 0.02651696237201191
$ws.Range("E10").Value = -0.01373232219717158
$ws.Range("D11").Value = 0.02425331924269382
$ws.Range("E11").Value = -0.002240896358543409
$ws.Range("D12").Value = 0.058436721859087
$ws.Range("E12").Value = 0.0002483238142538813
$ws.Range("D13").Value = 0.02634893646328202
$ws.Range("E13").Value = 0.005056843074501272
$ws.Range("D14").Value = 0.02777373268947532
$ws.Range("E14").Value = -0.01251604621309377
$ws.Range("D15").Value = 0.03581045386158876
$ws.Range("E15").Value = -0.004079551249362545
$ws.Range("D16").Value = 0.0193055877680095
$ws.Range("E16").Value = -0.01392681594756962
$ws.Range("D17").Value = 0.03024259830273832
$ws.Range("E17").Value = -0.001886961991194203
$ws.Range("D18").Value = 0.02389852782977213
$ws.Range("E18").Value = 0.003479471120389599
$ws.Range("D19").Value = 0.1327911213274999
$ws.Range("E19").Value = 0.001344086021505264
$ws.Range("D20").Value = 0.009787742884967459
$ws.Range("E20").Value = -0.01711366538952741
$ws.Range("D21").Value = 0.01605170267431745
$ws.Range("E21").Value = -0.008603526734925815
$ws.Range("D22").Value = 0.01730994041949976
$ws.Range("E22").Value = 0.0001067520683213363
$ws.Range("D23").Value = 0.01670584934189182
$ws.Range("E23").Value = 0.002127659574468144
$ws.Range("D24").Value = 0.02159167710650907
$ws.Range("E24").Value = -0.0007576577551683394
$ws.Range("D25").Value = 0.01194741604294376
$ws.Range("E25").Value = -0.009143553794574766
$ws.Range("D26").Value = 0.04407419153673129
$ws.Range("E26").Value = -0.007824205094056902
$ws.Range("D27").Value = 0.02548563243150914
$ws.Range("E27").Value = -0.00009809691975648516
$ws.Range("D28").Value = 0.04782200193076531
$ws.Range("E28").Value = 0.003631082062454549
$ws.Range("D29").Value = 0.0575125141420629
$ws.Range("E29").Value = -0.007408779403593191
$ws.Range("D30").Value = 0.01321717581317806
$ws.Range("E30").Value = 0.01979522184300331
$ws.Range("D31").Value = 0.01435497578706219
$ws.Range("E31").Value = 0.004732607666824373
$ws.Range("D32").Value = 0.04434713309244402
$ws.Range("E32").Value = -0.002073613271124986
$ws.Range("D33").Value = 0.01678824269084851
$ws.Range("E33").Value = 0.0009517766497462166
$ws.Range("E34").Value = -0.001607018139374139

# Restore sheet protection to its prior (protected) state.
$ws.Protect()
